# Auto-generated: update Sheets via scheduled runner
# Updates cached market-price / profit figures (columns H-N) across several leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 628.4167
$ws.Range("I12").Value = 621
$ws.Range("J12").Value = 635.8333
$ws.Range("K12").Value = 621
$ws.Range("L12").Value = 635.8333
$ws.Range("M12").Value = -451
$ws.Range("N12").Value = -975.8333
$ws.Range("H100").Value = 1666
$ws.Range("I100").Value = 999
$ws.Range("K100").Value = 999
$ws.Range("M100").Value = -458
$ws.Range("H107").Value = 939.7143
$ws.Range("I107").Value = 715
$ws.Range("K107").Value = 715
$ws.Range("M107").Value = 1205
$ws.Range("H112").Value = 3971.4285
$ws.Range("J112").Value = 3971.4285
$ws.Range("L112").Value = 11914.2855
$ws.Range("N112").Value = -14130.2855
$ws.Range("H132").Value = 1094.8125
$ws.Range("I132").Value = 1081.1333
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 3243.3999
$ws.Range("L132").Value = 3900
$ws.Range("M132").Value = -713.3998999999999
$ws.Range("N132").Value = -8960
$ws.Range("H138").Value = 2783.5134
$ws.Range("I138").Value = 3251.1428
$ws.Range("K138").Value = 9753.428400000001
$ws.Range("M138").Value = -4613.428400000001
$ws.Range("H141").Value = 4331.5557
$ws.Range("I141").Value = 3175.8
$ws.Range("K141").Value = 9527.400000000001
$ws.Range("M141").Value = -4347.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4074.6086
$ws.Range("I32").Value = 2746.4285
$ws.Range("K32").Value = 2746.4285
$ws.Range("M32").Value = -2459.4285
$ws.Range("H45").Value = 6430102.5
$ws.Range("I45").Value = 11251424
$ws.Range("J45").Value = 1673.8334
$ws.Range("K45").Value = 11251424
$ws.Range("L45").Value = 1673.8334
$ws.Range("M45").Value = -11251047
$ws.Range("N45").Value = -2427.8334
$ws.Range("H61").Value = 3089
$ws.Range("I61").Value = 2197.1765
$ws.Range("J61").Value = 6879.25
$ws.Range("K61").Value = 2197.1765
$ws.Range("L61").Value = 6879.25
$ws.Range("M61").Value = -1985.1765
$ws.Range("N61").Value = -7303.25
$ws.Range("H74").Value = 2793.4
$ws.Range("I74").Value = 1989.3334
$ws.Range("J74").Value = 3999.5
$ws.Range("K74").Value = 1989.3334
$ws.Range("L74").Value = 3999.5
$ws.Range("M74").Value = -1115.3334
$ws.Range("N74").Value = -5747.5
$ws.Range("H77").Value = 2793.4
$ws.Range("I77").Value = 1989.3334
$ws.Range("J77").Value = 3999.5
$ws.Range("K77").Value = 9946.666999999999
$ws.Range("L77").Value = 19997.5
$ws.Range("M77").Value = -5578.666999999999
$ws.Range("N77").Value = -28733.5
$ws.Range("H136").Value = 3089
$ws.Range("I136").Value = 2197.1765
$ws.Range("J136").Value = 6879.25
$ws.Range("K136").Value = 6591.529500000001
$ws.Range("L136").Value = 20637.75
$ws.Range("M136").Value = -4041.529500000001
$ws.Range("N136").Value = -25737.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 50000
$ws.Range("J28").Value = 50000
$ws.Range("L28").Value = 50000
$ws.Range("N28").Value = -50588
$ws.Range("H87").Value = 45000
$ws.Range("I87").Value = 40000
$ws.Range("J87").Value = 50000
$ws.Range("K87").Value = 40000
$ws.Range("L87").Value = 50000
$ws.Range("M87").Value = -38752
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 45000
$ws.Range("I90").Value = 40000
$ws.Range("J90").Value = 50000
$ws.Range("K90").Value = 120000
$ws.Range("L90").Value = 150000
$ws.Range("M90").Value = -113760
$ws.Range("N90").Value = -162480
$ws.Range("H94").Value = 714.4286
$ws.Range("I94").Value = 529.3
$ws.Range("K94").Value = 529.3
$ws.Range("M94").Value = -78.29999999999995
$ws.Range("H99").Value = 1399.4615
$ws.Range("I99").Value = 1299.5
$ws.Range("K99").Value = 1299.5
$ws.Range("M99").Value = 198.5
$ws.Range("H107").Value = 2199
$ws.Range("I107").Value = 1879.25
$ws.Range("J107").Value = 2838.5
$ws.Range("K107").Value = 1879.25
$ws.Range("L107").Value = 2838.5
$ws.Range("M107").Value = 40.75
$ws.Range("N107").Value = -6678.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 212.5
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 500
$ws.Range("N7").Value = -726
$ws.Range("H31").Value = 3769.2
$ws.Range("I31").Value = 1258.4
$ws.Range("J31").Value = 6280
$ws.Range("K31").Value = 1258.4
$ws.Range("L31").Value = 6280
$ws.Range("M31").Value = -963.4000000000001
$ws.Range("N31").Value = -6870
$ws.Range("H34").Value = 3769.2
$ws.Range("I34").Value = 1258.4
$ws.Range("J34").Value = 6280
$ws.Range("K34").Value = 1258.4
$ws.Range("L34").Value = 6280
$ws.Range("M34").Value = -1056.4
$ws.Range("N34").Value = -6684
$ws.Range("H94").Value = 1380.6
$ws.Range("I94").Value = 1322.6
$ws.Range("J94").Value = 1438.6
$ws.Range("K94").Value = 1322.6
$ws.Range("L94").Value = 1438.6
$ws.Range("M94").Value = -871.5999999999999
$ws.Range("N94").Value = -2340.6
$ws.Range("H107").Value = 455.6
$ws.Range("I107").Value = 362.16666
$ws.Range("J107").Value = 829.3333
$ws.Range("K107").Value = 362.16666
$ws.Range("L107").Value = 829.3333
$ws.Range("M107").Value = 1557.83334
$ws.Range("N107").Value = -4669.3333
$ws.Range("H132").Value = 2846.4707
$ws.Range("I132").Value = 1172.5
$ws.Range("K132").Value = 3517.5
$ws.Range("M132").Value = -987.5
$ws.Range("H134").Value = 2867.4167
$ws.Range("I134").Value = 2512.2222
$ws.Range("K134").Value = 7536.6666
$ws.Range("M134").Value = -5001.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11127673
$ws.Range("J131").Value = 20015.676
$ws.Range("L131").Value = 60047.028
$ws.Range("N131").Value = -70127.02799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3644464.8
$ws.Range("J24").Value = 12730.143
$ws.Range("L24").Value = 12730.143
$ws.Range("N24").Value = -13076.143
$ws.Range("H132").Value = 3739.6
$ws.Range("J132").Value = 4126.2856
$ws.Range("L132").Value = 12378.8568
$ws.Range("N132").Value = -17438.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4068.476
$ws.Range("I7").Value = 1951.6154
$ws.Range("J7").Value = 7508.375
$ws.Range("K7").Value = 1951.6154
$ws.Range("L7").Value = 7508.375
$ws.Range("M7").Value = -1839.6154
$ws.Range("N7").Value = -7732.375
$ws.Range("H40").Value = 7285.7856
$ws.Range("J40").Value = 11085.571
$ws.Range("L40").Value = 11085.571
$ws.Range("N40").Value = -11357.571
$ws.Range("H46").Value = 1353.6923
$ws.Range("I46").Value = 399.5
$ws.Range("K46").Value = 399.5
$ws.Range("M46").Value = -211.5
$ws.Range("H61").Value = 2182.5557
$ws.Range("I61").Value = 2078.7334
$ws.Range("J61").Value = 2701.6667
$ws.Range("K61").Value = 2078.7334
$ws.Range("L61").Value = 2701.6667
$ws.Range("M61").Value = -1876.7334
$ws.Range("N61").Value = -3105.6667
$ws.Range("H113").Value = 2182.5557
$ws.Range("I113").Value = 2078.7334
$ws.Range("J113").Value = 2701.6667
$ws.Range("K113").Value = 2078.7334
$ws.Range("L113").Value = 2701.6667
$ws.Range("M113").Value = 91.26659999999993
$ws.Range("N113").Value = -7041.6667
$ws.Range("H126").Value = 4068.476
$ws.Range("I126").Value = 1951.6154
$ws.Range("J126").Value = 7508.375
$ws.Range("K126").Value = 5854.8462
$ws.Range("L126").Value = 22525.125
$ws.Range("M126").Value = -3384.8462
$ws.Range("N126").Value = -27465.125
$ws.Range("H132").Value = 1797.76
$ws.Range("I132").Value = 1416.875
$ws.Range("J132").Value = 1977
$ws.Range("K132").Value = 4250.625
$ws.Range("L132").Value = 5931
$ws.Range("M132").Value = -1720.625
$ws.Range("N132").Value = -10991

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 47454.582
$ws.Range("J123").Value = 47454.582
$ws.Range("L123").Value = 47454.582
$ws.Range("N123").Value = -57254.582
$ws.Range("H126").Value = 5322.7393
$ws.Range("J126").Value = 7467.1665
$ws.Range("L126").Value = 22401.4995
$ws.Range("N126").Value = -27341.4995
$ws.Range("H132").Value = 1777.1428
$ws.Range("I132").Value = 908.13336
$ws.Range("J132").Value = 3949.6667
$ws.Range("K132").Value = 2724.40008
$ws.Range("L132").Value = 11849.0001
$ws.Range("M132").Value = -194.4000800000003
$ws.Range("N132").Value = -16909.0001
$ws.Range("H136").Value = 2069.3696
$ws.Range("I136").Value = 1775.4375
$ws.Range("J136").Value = 2741.2144
$ws.Range("K136").Value = 5326.3125
$ws.Range("L136").Value = 8223.643199999999
$ws.Range("M136").Value = -2776.3125
$ws.Range("N136").Value = -13323.6432

Write-Output "Applied 223 cell updates across 8 sheets"